$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Text
    )
    $cell = $ws.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.ClearFormats()
}

Set-TextValue "D2" '28.991.34'
Set-TextValue "E2" '  +0.48%  '
Set-TextValue "D3" '1.886.69'
Set-TextValue "E3" '  -0.95%  '
Set-TextValue "D4" '1.000'
Set-TextValue "E4" '  +0.01%  '
Set-TextValue "D5" '330.41'
Set-TextValue "E5" '  -2.32%  '
Set-TextValue "E6" '  +0.06%  '
Set-TextValue "D7" '0.4592'
Set-TextValue "E7" '  -2.54%  '
Set-TextValue "D8" '0.4060'
Set-TextValue "E8" '  +0.61%  '
Set-TextValue "D9" '47.63'
Set-TextValue "E9" '  -0.58%  '
Set-TextValue "D10" '0.07970'
Set-TextValue "E10" '  -2.24%  '
Set-TextValue "D11" '0.9923'
Set-TextValue "E11" '  -2.47%  '
Set-TextValue "D12" '21.66'
Set-TextValue "E12" '  -3.21%  '
Set-TextValue "D13" '1.878.39'
Set-TextValue "E13" '  -0.55%  '
Set-TextValue "D14" '5.907'
Set-TextValue "E14" '  -2.90%  '
Set-TextValue "D15" '7.064'
Set-TextValue "E15" '  -3.76%  '
Set-TextValue "E16" '  -0.01%  '
Set-TextValue "D17" '88.40'
Set-TextValue "E17" '  -3.09%  '
Set-TextValue "D18" '0.00001026'
Set-TextValue "E18" '  -2.45%  '
Set-TextValue "D19" '0.06551'
Set-TextValue "E19" '  -0.91%  '
Set-TextValue "D20" '17.41'
Set-TextValue "E20" '  -2.02%  '
Set-TextValue "E21" '  -0.04%  '
Set-TextValue "D22" '28.989.63'
Set-TextValue "E22" '  +0.36%  '
Set-TextValue "E23" '  -2.48%  '
Set-TextValue "D24" '11.29'
Set-TextValue "E24" '  +1.43%  '
Set-TextValue "D25" '2.209'
Set-TextValue "E25" '  -2.59%  '
Set-TextValue "D26" '2.106.02'
Set-TextValue "E26" '  -0.39%  '
Set-TextValue "D27" '156.63'
Set-TextValue "E27" '  -2.86%  '
Set-TextValue "D28" '19.59'
Set-TextValue "E28" '  -2.01%  '
Set-TextValue "D29" '2.103'
Set-TextValue "E29" '  -2.75%  '
Set-TextValue "D30" '5.407'
Set-TextValue "E30" '  -2.21%  '
Set-TextValue "D31" '117.80'
Set-TextValue "E31" '  -2.14%  '
Set-TextValue "E32" '  -1.40%  '
Set-TextValue "D33" '0.09330'
Set-TextValue "E33" '  -2.50%  '
Set-TextValue "D34" '3.598'
Set-TextValue "E34" '  -1.57%  '
Set-TextValue "D35" '1.406'
Set-TextValue "E35" '  +0.17%  '
Set-TextValue "D36" '5.281'
Set-TextValue "E36" '  -2.19%  '
Set-TextValue "D37" '0.06049'
Set-TextValue "E37" '  -2.41%  '
Set-TextValue "E38" '  -3.06%  '
Set-TextValue "D39" '8.269'
Set-TextValue "E39" '  -4.45%  '
Set-TextValue "D40" '1.180'
Set-TextValue "E40" '  -1.22%  '
Set-TextValue "D41" '1.000'
Set-TextValue "E41" '  +0.04%  '
Set-TextValue "E42" '  -3.82%  '
Set-TextValue "B43" 'Aptos'
Set-TextValue "C43" 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue "D43" '10.12'
Set-TextValue "E43" '  -3.38%  '
Set-TextValue "B44" 'Algorand'
Set-TextValue "C44" 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue "D44" '0.1819'
Set-TextValue "E44" '  -4.26%  '
Set-TextValue "D45" '1.260'
Set-TextValue "E45" '  -0.63%  '
Set-TextValue "D46" '0.07525'
Set-TextValue "E46" '  +3.72%  '
Set-TextValue "D47" '2.270'
Set-TextValue "E47" '  +6.97%  '
Set-TextValue "D48" '12.03'
Set-TextValue "E48" '  -2.35%  '
Set-TextValue "E49" '  -3.05%  '
Set-TextValue "E50" '  -4.16%  '
Set-TextValue "B51" 'Elrond'
Set-TextValue "C51" 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextValue "D51" '45.65'
Set-TextValue "E51" '  +12.78%  '
